$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.372.29"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.23%  '
$ws.Range('D3').Value = "'1.857.86"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -4.03%  '
$ws.Range('E4').Value = '  -1.10%  '
$ws.Range('D5').Value = "'323.59"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.77%  '
$ws.Range('D6').Value = "'1.002"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.93%  '
$ws.Range('E7').Value = '  -4.78%  '
$ws.Range('D8').Value = "'0.3867"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.72%  '
$ws.Range('D9').Value = "'48.89"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -8.72%  '
$ws.Range('D10').Value = "'0.07904"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.11%  '
$ws.Range('E11').Value = '  -3.35%  '
$ws.Range('E12').Value = '  -4.29%  '
$ws.Range('D13').Value = "'1.854.93"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.03%  '
$ws.Range('D14').Value = "'5.924"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.35%  '
$ws.Range('D15').Value = "'7.116"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.52%  '
$ws.Range('D16').Value = "'1.002"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.20%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').Value = "'85.90"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.59%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = "'0.00001032"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.65%  '
$ws.Range('D19').Value = "'0.06519"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.47%  '
$ws.Range('D20').Value = "'17.04"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.54%  '
$ws.Range('E21').Value = '  -1.08%  '
$ws.Range('D22').Value = "'5.527"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.77%  '
$ws.Range('D23').Value = "'27.373.99"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E24').Value = '  -5.41%  '
$ws.Range('D25').Value = "'2.284"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.73%  '
$ws.Range('D26').Value = "'2.075.78"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.14%  '
$ws.Range('D27').Value = "'153.83"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.00%  '
$ws.Range('D28').Value = "'19.78"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.22%  '
$ws.Range('D29').Value = "'2.075"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.60%  '
$ws.Range('D30').Value = "'5.439"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.00%  '
$ws.Range('D31').Value = "'121.01"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.46%  '
$ws.Range('D32').Value = "'1.483"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.50%  '
$ws.Range('D33').Value = "'0.09285"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.49%  '
$ws.Range('D34').Value = "'0.9370"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.97%  '
$ws.Range('D35').Value = "'3.601"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.73%  '
$ws.Range('D36').Value = "'5.251"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.25%  '
$ws.Range('D37').Value = "'0.02236"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.74%  '
$ws.Range('D38').Value = "'1.222"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.66%  '
$ws.Range('D39').Value = "'0.05986"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.20%  '
$ws.Range('D40').Value = "'8.208"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -11.54%  '
$ws.Range('E41').Value = '  -0.98%  '
$ws.Range('D42').Value = "'0.5907"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.82%  '
$ws.Range('D43').Value = "'0.1892"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.06%  '
$ws.Range('D44').Value = "'10.09"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -9.64%  '
$ws.Range('D45').Value = "'1.278"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.45%  '
$ws.Range('D46').Value = "'0.5609"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.36%  '
$ws.Range('D47').Value = "'11.97"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.94%  '
$ws.Range('D48').Value = "'3.364"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.97%  '
$ws.Range('D49').Value = "'1.922"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.23%  '
$ws.Range('E50').Value = '  -0.31%  '
$ws.Range('D51').Value = "'108.22"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.72%  '
